$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("BD1").Value = 0.94160910907091488
$ws.Range("G2").Value = 0.8918706253105958
$ws.Range("J2").Value = 0.60447031948537189
$ws.Range("Z2").Value = 0.96136470050596479
$ws.Range("T3").Value = 0.78643969666329716
$ws.Range("AR3").Value = 0.72876410708156558
$ws.Range("BG4").Value = 0.90232724304789214
$ws.Range("AB5").Value = 0.75104594943086056
$ws.Range("BJ5").Value = 0.92515578787683928
$ws.Range("BN5").Value = 0.89441174562515502
$ws.Range("D6").Value = 0.89884817793018734
$ws.Range("H6").Value = 0.54961256774389433
$ws.Range("N7").Value = 0.7854191774931385
$ws.Range("U8").Value = 0.77371296911921505
$ws.Range("AF8").Value = 0.98951079063530356
$ws.Range("J9").Value = 0.64947607059273205
$ws.Range("H10").Value = 0.98651649705565481
$ws.Range("M10").Value = 0.95574922254489414
$ws.Range("AN10").Value = 0.55549008901567021
$ws.Range("E11").Value = 0.90712771128691461
$ws.Range("M11").Value = 0.96247716473639588
$ws.Range("T11").Value = 0.65814780328745592
$ws.Range("AN11").Value = 0.73225481253602309
$ws.Range("AB12").Value = 0.79802189769693643
$ws.Range("AF12").Value = 0.67712480660732965
$ws.Range("BA12").Value = 0.87334513469538344
$ws.Range("L14").Value = 0.97615082274485454
$ws.Range("AP14").Value = 0.91001011943217103
$ws.Range("AW14").Value = 0.95896482911635506
$ws.Range("J15").Value = 0.94629642625831512
$ws.Range("AI15").Value = 0.91985984183505654
$ws.Range("B16").Value = 0.81381375503094611
$ws.Range("AD16").Value = 0.99303027845108627
$ws.Range("AK16").Value = 0.78974124690948844
$ws.Range("AV16").Value = 0.85176051332411618
$ws.Range("R17").Value = 0.872737762470704
$ws.Range("AI17").Value = 0.95424522003725332
$ws.Range("BB18").Value = 0.99467643609059708
$ws.Range("J19").Value = 0.75658606172635479
$ws.Range("S20").Value = 0.96654067142056133
$ws.Range("V21").Value = 0.95363820547467348
$ws.Range("AD21").Value = 0.93586045296907416
$ws.Range("BG22").Value = 0.7458162869179058
$ws.Range("BK22").Value = 0.58010344298816441
$ws.Range("Y23").Value = 0.88894867369528807
$ws.Range("AT23").Value = 0.99383149033948825
$ws.Range("BG24").Value = 0.83622921137293216
$ws.Range("BH24").Value = 0.885663395929283
$ws.Range("BL24").Value = 0.68771738418691974
$ws.Range("AA25").Value = 0.88949041804757356
$ws.Range("V26").Value = 0.6897389118373034
$ws.Range("AL26").Value = 0.84417963750213976
$ws.Range("AN26").Value = 0.96979597947275975
$ws.Range("AO26").Value = 0.73182057902712283
$ws.Range("AQ27").Value = 0.87305809774559306
$ws.Range("W29").Value = 0.89686689471646852
$ws.Range("AB29").Value = 0.98496007230729177
$ws.Range("AF29").Value = 0.66010892562030843
$ws.Range("AL29").Value = 0.68604413143276255
$ws.Range("M30").Value = 0.98347315950100622
$ws.Range("AC31").Value = 0.60754227565767693
$ws.Range("E32").Value = 0.87789728439223702
$ws.Range("AE32").Value = 0.74450038938586427
$ws.Range("AB33").Value = 0.69545573848408593
$ws.Range("V35").Value = 0.95651393240036209
$ws.Range("AG35").Value = 0.98447102421434785
$ws.Range("BN35").Value = 0.8574481965597659
$ws.Range("AY36").Value = 0.91153956028221617
$ws.Range("I37").Value = 0.7277859963885176
$ws.Range("AJ37").Value = 0.73278653544394112
$ws.Range("AM37").Value = 0.64300276123202127
$ws.Range("AU37").Value = 0.88962027921348696
$ws.Range("K38").Value = 0.90647539161615065
$ws.Range("M38").Value = 0.9692349537133006
$ws.Range("R38").Value = 0.97822271784803849
$ws.Range("Q39").Value = 0.8712355275756174
$ws.Range("AY39").Value = 0.66232057514116338
$ws.Range("AZ41").Value = 0.94606861854944735
$ws.Range("BE41").Value = 0.95036160721282692
$ws.Range("R42").Value = 0.96754809151858756
$ws.Range("AK42").Value = 0.8592822348182636
$ws.Range("T43").Value = 0.91153976192644293
$ws.Range("BH44").Value = 0.84971236607136125
$ws.Range("X45").Value = 0.71690711690181486
$ws.Range("AR45").Value = 0.85142180485746755
$ws.Range("AT45").Value = 0.9221421869028863
$ws.Range("M46").Value = 0.97409004262110876
$ws.Range("O47").Value = 0.94669887284085552
$ws.Range("Z47").Value = 0.99651134967731347
$ws.Range("AT47").Value = 0.69613632841491024
$ws.Range("AV47").Value = 0.90452056594349672
$ws.Range("AP48").Value = 0.86041581901356534
$ws.Range("BH49").Value = 0.99695805251396874
$ws.Range("AZ50").Value = 0.78290925694882874
$ws.Range("BD50").Value = 0.77455184622606343
$ws.Range("BM50").Value = 0.75668548651524592
$ws.Range("AQ51").Value = 0.76961563926178389
$ws.Range("AW51").Value = 0.77001466137576868
$ws.Range("J52").Value = 0.91517920367161865
$ws.Range("X52").Value = 0.97499779881999837
$ws.Range("F53").Value = 0.99084508530264848
$ws.Range("AM53").Value = 0.92317479909683842
$ws.Range("AQ53").Value = 0.97464162493339157
$ws.Range("BK53").Value = 0.56644714773132665
$ws.Range("BA54").Value = 0.79401533006437686
$ws.Range("BJ54").Value = 0.96767829739582567
$ws.Range("AO55").Value = 0.76726787203849722
$ws.Range("AU55").Value = 0.66433907790112556
$ws.Range("BA55").Value = 0.75922381844817899
$ws.Range("I56").Value = 0.89846744654992583
$ws.Range("AE56").Value = 0.95738274307690285
$ws.Range("A57").Value = 0.86157917753228963
$ws.Range("AQ57").Value = 0.91250483626278756
$ws.Range("BG57").Value = 0.99628247966761974
$ws.Range("R58").Value = 0.69130091192895993
$ws.Range("J59").Value = 0.78095564142310547
$ws.Range("AA59").Value = 0.77107106041522178
$ws.Range("AP59").Value = 0.87518567579211015
$ws.Range("Y60").Value = 0.81025846494910569
$ws.Range("AH60").Value = 0.96838190407223257
$ws.Range("X62").Value = 0.95798527786121634
$ws.Range("AE62").Value = 0.94125874806947585
$ws.Range("BF62").Value = 0.85083235114829581
$ws.Range("BH62").Value = 0.82352187495538121
$ws.Range("BI62").Value = 0.84866452567017781
$ws.Range("P63").Value = 0.99357680880940391
$ws.Range("BI63").Value = 0.98110862002137578
$ws.Range("AR64").Value = 0.56048948047259151
$ws.Range("BF64").Value = 0.72132317055924666
$ws.Range("BO64").Value = 0.64863386786632127
$ws.Range("AH65").Value = 0.93911489202644782
$ws.Range("BP65").Value = 0.81741172196899947
$ws.Range("AO66").Value = 0.81991274694620997
$ws.Range("P67").Value = 0.91605470970018865
$ws.Range("B68").Value = 0.85781688885492935
$ws.Range("AV68").Value = 0.89040892103103486
